# Generate Report for Handoff
# "b.md" has now been handed off again (new .xlf handback files produced),
# so its status flips from "in sync" to "ready for handoff" and its
# handback file / datetime / error-detail get refreshed on every sheet.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/230af4162d4dedce5deed2aaecd6e25474b45204/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/506cb60debfde8aaff1b4ac2a7a3376a24db09ef/e2e/b.md."

# --- Overview sheet: row 3 is the b.md entry -----------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-24 22:38:14"

# --- zh-cn sheet: row 3 is the b.md entry --------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "False"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-24 22:38:08"
$zh.Range("P3").Value = $errorDetail
$zh.Columns.Item(16).ColumnWidth = 40

# --- de-de sheet: row 3 is the b.md entry --------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "False"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-08-24 22:38:14"
$de.Range("P3").Value = $errorDetail
$de.Columns.Item(16).ColumnWidth = 40
